$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.391.98'
$ws.Range('E2').Value = '  +1.61%  '

# Row 3
$ws.Range('D3').Value = '1.955.61'
$ws.Range('E3').Value = '  +3.44%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '328.00'
$ws.Range('E5').Value = '  +0.75%  '

# Row 6
$ws.Range('E6').Value = '  +0.13%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4642'
$ws.Range('E7').Value = '  +1.48%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3935'
$ws.Range('E8').Value = '  +0.87%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '46.30'
$ws.Range('E9').Value = '  -0.79%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07885'
$ws.Range('E10').Value = '  +0.62%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9987'
$ws.Range('E11').Value = '  +1.24%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.34'
$ws.Range('E12').Value = '  +2.13%  '

# Row 13
$ws.Range('D13').Value = '1.986.57'
$ws.Range('E13').Value = '  +5.93%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.149'
$ws.Range('E14').Value = '  +1.58%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.828'
$ws.Range('E15').Value = '  +2.56%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.07124'
$ws.Range('E16').Value = '  +2.83%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '88.49'
$ws.Range('E17').Value = '  +0.58%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.006'
$ws.Range('E18').Value = '  +0.31%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000009940'
$ws.Range('E19').Value = '  -0.15%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.09'
$ws.Range('E20').Value = '  +0.73%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.09%  '

# Row 22
$ws.Range('D22').Value = '29.426.40'
$ws.Range('E22').Value = '  +1.72%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.527'
$ws.Range('E23').Value = '  +4.41%  '

# Row 24
$ws.Range('E24').Value = '  +2.42%  '

# Row 25
$ws.Range('D25').Value = '2.206.18'
$ws.Range('E25').Value = '  +1.82%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.124'
$ws.Range('E26').Value = '  +3.39%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '158.37'
$ws.Range('E27').Value = '  +1.64%  '

# Row 28
$ws.Range('E28').Value = '  +1.51%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.965'
$ws.Range('E29').Value = '  +1.74%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '120.02'
$ws.Range('E30').Value = '  +2.20%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.877'
$ws.Range('E31').Value = '  -2.56%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09377'
$ws.Range('E32').Value = '  +0.50%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8927'
$ws.Range('E33').Value = '  -1.30%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.250'
$ws.Range('E34').Value = '  -0.72%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.335'
$ws.Range('E35').Value = '  +1.12%  '

# Row 36
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.174'
$ws.Range('E36').Value = '  -2.74%  '

# Row 37
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.000003741'
$ws.Range('E37').Value = '  +127.37%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05809'
$ws.Range('E38').Value = '  +0.72%  '

# Row 39
$ws.Range('E39').Value = '  -1.03%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02116'
$ws.Range('E40').Value = '  +2.21%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.002'
$ws.Range('E41').Value = '  +0.09%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '7.806'
$ws.Range('E42').Value = '  +1.79%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5750'
$ws.Range('E43').Value = '  +1.40%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1822'
$ws.Range('E44').Value = '  +3.14%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '9.809'
$ws.Range('E45').Value = '  +0.56%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.07'
$ws.Range('E46').Value = '  +1.04%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5366'
$ws.Range('E47').Value = '  +0.35%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.190'
$ws.Range('E48').Value = '  -4.15%  '

# Row 49
$ws.Range('E49').Value = '  +1.52%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06953'
$ws.Range('E50').Value = '  -1.27%  '

# Row 51
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.605'
$ws.Range('E51').Value = '  +3.33%  '
